## Add files via upload
## The workbook gains a new data row (row 50) that is effectively a
## duplicate of row 49 with a handful of cells edited, plus two brand new
## shared strings introduced by those edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 49 (values + formatting) down into the new row 50.
$ws.Rows.Item(49).Copy()
$ws.Rows.Item(50).Insert()

# Match the (wrapped-text driven) row height of the source row.
$ws.Rows.Item(50).RowHeight = 236.25

# Now apply the specific cell edits that differentiate row 50 from row 49.
# NOTE: the two brand new shared strings must be introduced in the same
# order they appear in the workbook's shared string table, so set the
# "Variable Description" (S50) text before the "Activity Name" (E50) text.
$ws.Range("D50").Value = 1.2
$ws.Range("S50").Value = "Chlorine Concentration in Cooling Water 11"
$ws.Range("E50").Value = "New Activity Test  19"
$ws.Range("R50").Value = "Variable2"

# Reflect where the user ended up after making the edit: scrolled down so
# row 49 is the first visible row, with G50 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G50").Select()
